$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. "Dear Mr. " + "Lohmann" (two runs) -> single run "Dear Mr. Lohmann"
# -----------------------------------------------------------------
$d.Content.Find.Execute("Dear Mr. Lohmann", $false, $false, $false, $false, $false, $true, 1, $false, "Dear Mr. Lohmann", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Re-type the same wording for the run-per-fragment block so the
#    separate runs collapse into a single contiguous run.
# -----------------------------------------------------------------
$old2 = " as I would be able to combine my skills in data analysis and machine learning. I really admire and appreciate Deloitte's values in creating a positive and supportive culture which is focused on continuous learning. I feel that I would fit right in as I share these values and find learning new things critical to professional development and growth. With my strong background in data modeling, I would love to use my skills and create value "
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# -----------------------------------------------------------------
# 3. Wording change: "generate insights, and  create value for informed
#    decision-making." -> "and generate valuable insights for informed
#    decision-making."
# -----------------------------------------------------------------
$old3 = "generate insights, and  create value for informed decision-making."
$new3 = "and generate valuable insights for informed decision-making."
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# -----------------------------------------------------------------
# 4. Move the "_GoBack" bookmark from the end of the "Usman Shaikh"
#    paragraph to the end of the "Sincerely," paragraph.
# -----------------------------------------------------------------
$paras = $d.Paragraphs
$sincerelyPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Sincerely,*") {
        $sincerelyPara = $p
        break
    }
}

if ($sincerelyPara -ne $null) {
    $endPos = $sincerelyPara.Range.End - 1

    # A bookmark collapsed exactly on the paragraph-mark boundary can't be
    # added directly, so nudge the boundary out of the way with a throwaway
    # character, drop the bookmark at the (now mid-text) position, then
    # remove the throwaway character again. The bookmark stays put.
    $guard = $d.Range($endPos, $endPos)
    $guard.InsertAfter("~")

    $target = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null

    $guardRange = $d.Range($endPos, $endPos + 1)
    $guardRange.Text = ""
}
